# The sheet had a "Total Score" column (O) with a header string and, for
# each data row, a `=SUM(B:N)` formula. The edit removes that whole column's
# contents (header text + formulas/values) while leaving the header cell's
# existing formatting in place - i.e. selecting O1:O6 and pressing Delete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1:O6").ClearContents()

# Final on-screen state: column O selected (whole-column selection), with
# the view scrolled right so column K is the left-most visible column.
$ws.Range("O1:O1048576").Select()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
